$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Service Code"
$ws.Range("B1").Value = "Service Title"
$ws.Range("C1").Value = "Service Duration (minutes)"
$ws.Range("D1").ClearContents()

# D1 used to carry the bold header style (s=1); the new layout drops the
# 4th header but keeps a normal-bodied, wrap-text style (s=2) on it instead.
# Copy the body style from A2 (already style index 2) onto D1 without
# touching D1's (now empty) value.
$ws.Range("A2").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 1 grew taller to fit the new, longer header text.
$ws.Rows("1").RowHeight = 57.6

# --- Data rows (rows 2-4): Service Code / Service Title / Duration ---
$ws.Range("A2").Value = "A"
$ws.Range("B2").Value = "Consultation"
$ws.Range("C2").Value = 5
$ws.Range("D2").ClearContents()

$ws.Range("A3").Value = "B"
$ws.Range("B3").Value = "Follow-up"
$ws.Range("C3").Value = 3
$ws.Range("D3").ClearContents()

$ws.Range("A4").Value = "C"
$ws.Range("B4").Value = "Diagnostic"
$ws.Range("C4").Value = 7
$ws.Range("D4").ClearContents()

# --- Rows 5-7 are now fully empty (previously held extra sample rows) ---
$ws.Range("A5:D7").ClearContents()

# --- Selection moves to I2 ---
$ws.Range("I2").Select()
